$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for rows 2 through 9
# from serial date 46075 (2026-02-22) to 46076 (2026-02-23).
foreach ($row in 2..9) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value2 = 46076
    }
}
